$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.336.95"
$ws.Range("E2").Value = "  +0.39%  "
$ws.Range("D3").Value = "3.741.39"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.95%  "
$ws.Range("D7").Value = "3.737.12"
$ws.Range("E7").Value = "  -1.08%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.518"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -3.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.39"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000257"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.28%  "
$ws.Range("D15").Value = "4.372.30"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "3.747.82"
$ws.Range("E16").Value = "  -0.77%  "
$ws.Range("D17").Value = "68.283.25"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.49%  "
$ws.Range("E19").Value = "  +0.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "463.16"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.690"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.84"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000143"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.54%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.32%  "
$ws.Range("D30").Value = "3.893.06"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.07%  "
$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  --%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").Value = "3.698.00"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.02%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.137"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.22%  "
$ws.Range("B41").Value = "Mantle"
$ws.Range("C41").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.996"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.299"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.24"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +8.13%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.18%  "
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "145.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "384.79"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.38%  "
